# Updated symbol list on Mon Dec 26 19:16:00 UTC 2022 with GitHub Actions
# Refresh the "Price" column (D) for the rows whose coin price changed.
# NumberFormat is forced to Text ("@") before assignment so the numeric-looking
# strings are stored verbatim (matching significant digits / trailing zeros)
# instead of being auto-coerced into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2"  "242.73"
Set-TextValue "D3"  "23.10"
Set-TextValue "D4"  "5.415"
Set-TextValue "D5"  "0.05884"
Set-TextValue "D6"  "3.436"
Set-TextValue "D7"  "6.556"
Set-TextValue "D9"  "0.9409"
Set-TextValue "D10" "0.1418"
Set-TextValue "D11" "0.07427"
Set-TextValue "D12" "0.03360"
Set-TextValue "D13" "0.03056"
Set-TextValue "D14" "0.09334"
Set-TextValue "D16" "0.001572"
Set-TextValue "D17" "0.04666"
Set-TextValue "D18" "0.0005925"
Set-TextValue "D19" "0.005889"
Set-TextValue "D20" "0.001268"
Set-TextValue "D21" "0.004894"
Set-TextValue "D23" "3.565"
Set-TextValue "D27" "0.0002286"
Set-TextValue "D40" "0.03964"
Set-TextValue "D41" "0.006184"
Set-TextValue "D42" "0.1069"
Set-TextValue "D44" "0.009091"
Set-TextValue "D45" "0.00005200"
Set-TextValue "D47" "0.6706"
Set-TextValue "D48" "0.002387"
